# Update "想去人数" (F column) figures on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - rows keyed by their F-column (想去人数) values.
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 543   # 南宁·布谷鸟动漫展5th: 538 -> 543
$ws1.Range("F3").Value = 6476  # 南宁·2024良牙动漫秋季盛典（秋典）: 6426 -> 6476
$ws1.Range("F4").Value = 398   # 南宁·快看漫画动漫游戏嘉年华 KKWORLD-mini（取消）: 397 -> 398
$ws1.Range("F5").Value = 98    # 广西·THO04-永夜廻想: 97 -> 98
$ws1.Range("F6").Value = 140   # 南宁·花海演绎二次元水上派对: 137 -> 140
$ws1.Range("F8").Value = 75    # 南宁·熊喵M动漫嘉年华·万圣派对: 74 -> 75
$ws1.Range("F9").Value = 575   # 南宁·万圣漫控嘉年华10: 572 -> 575

# Sheet "全部类型" (All Types) - same events, different row offsets because of
# an extra performance row inserted earlier in this sheet.
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 543   # 南宁·布谷鸟动漫展5th: 538 -> 543
$ws4.Range("F3").Value = 6476  # 南宁·2024良牙动漫秋季盛典（秋典）: 6426 -> 6476
$ws4.Range("F4").Value = 398   # 南宁·快看漫画动漫游戏嘉年华 KKWORLD-mini（取消）: 397 -> 398
$ws4.Range("F6").Value = 98    # 广西·THO04-永夜廻想: 97 -> 98
$ws4.Range("F7").Value = 140   # 南宁·花海演绎二次元水上派对: 137 -> 140
$ws4.Range("F10").Value = 75   # 南宁·熊喵M动漫嘉年华·万圣派对: 74 -> 75
$ws4.Range("F11").Value = 575  # 南宁·万圣漫控嘉年华10: 572 -> 575
